# Update stats for 2025-08 (row 21) in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B21").Value = 6246
$ws.Range("C21").Value = 992
$ws.Range("D21").Value = 5646314
$ws.Range("E21").Value = 903.9887928274095
$ws.Range("F21").Value = 8.418677312966505
$ws.Range("G21").Value = 4.641350210970474
$ws.Range("H21").Value = 28.88267717757611
